# Issue #12 - "Simulate the brain" slide added on slides
#
# A new slide discussing how a human sleeps to organize memory (and why SAI
# needs something similar) is inserted right after the "issues / Hard to
# deploy" slide and right before the "Solutions" slide, i.e. it becomes the
# 8th slide of the deck (everything from the old "Solutions" slide onward
# shifts down by one position).

$p = $ppt.ActivePresentation

# Slide 8 is currently "Solutions" - reuse its layout ("Title and Content")
# for the newly inserted slide so it keeps the same look and feel.
$refSlide = $p.Slides.Item(8)
$layout = $refSlide.CustomLayout

$newSlide = $p.Slides.AddSlide(8, $layout)

$title = $newSlide.Shapes.Item(1)
$title.Name = "Titre 1"
$title.TextFrame.TextRange.Text = "Simulate the brain"

$body = $newSlide.Shapes.Item(2)
$body.Name = "Espace réservé du contenu 2"
$body.TextFrame.TextRange.Text = "A human sleep to organize it memory`rSo do SAI need to`rSchema to explain"
